# Remove trailing newline characters from the "Market Cap", "SEDOL" and
# "ISIN" header cells (F1, G1, H1) on every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("F1").Value = "Market Cap"
    $ws.Range("G1").Value = "SEDOL"
    $ws.Range("H1").Value = "ISIN"
}
